$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'91.595.01"
$ws.Range("E2").Value = "  +4.00%  "
$ws.Range("D3").Value = "'3.105.75"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'219.17"
$ws.Range("E5").Value = "  +4.63%  "
$ws.Range("D6").Value = "'621.32"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "'0.383"
$ws.Range("E7").Value = "  +3.64%  "
$ws.Range("D8").Value = "'0.894"
$ws.Range("E8").Value = "  +9.32%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "'3.101.83"
$ws.Range("E10").Value = "  +1.28%  "
$ws.Range("D11").Value = "'0.725"
$ws.Range("E11").Value = "  +22.47%  "
$ws.Range("D12").Value = "'0.190"
$ws.Range("E12").Value = "  +7.10%  "
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  +7.73%  "
$ws.Range("D14").Value = "'5.42"
$ws.Range("E14").Value = "  +3.18%  "
$ws.Range("D15").Value = "'91.173.68"
$ws.Range("E15").Value = "  +3.79%  "
$ws.Range("D16").Value = "'33.48"
$ws.Range("E16").Value = "  +6.69%  "
$ws.Range("D17").Value = "'3.680.22"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "'3.107.72"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "'3.74"
$ws.Range("E19").Value = "  +15.75%  "
$ws.Range("D20").Value = "'0.0000234"
$ws.Range("E20").Value = "  +11.97%  "
$ws.Range("D21").Value = "'13.86"
$ws.Range("E21").Value = "  +6.43%  "
$ws.Range("D22").Value = "'433.34"
$ws.Range("E22").Value = "  +3.41%  "
$ws.Range("D23").Value = "'8.68"
$ws.Range("E23").Value = "  +6.30%  "
$ws.Range("D24").Value = "'5.16"
$ws.Range("E24").Value = "  +7.39%  "
$ws.Range("D25").Value = "'5.56"
$ws.Range("E25").Value = "  +2.88%  "
$ws.Range("D26").Value = "'12.13"
$ws.Range("E26").Value = "  +7.92%  "
$ws.Range("D27").Value = "'84.38"
$ws.Range("E27").Value = "  +3.39%  "
$ws.Range("D28").Value = "'3.264.62"
$ws.Range("E28").Value = "  +1.25%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").Value = "'0.171"
$ws.Range("E30").Value = "  +17.61%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").Value = "'0.993"
$ws.Range("E31").Value = "  -8.65%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'8.77"
$ws.Range("E32").Value = "  +9.88%  "
$ws.Range("B33").Value = "dogwifhat"
$ws.Range("C33").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D33").Value = "'3.88"
$ws.Range("E33").Value = "  +10.34%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "'520.37"
$ws.Range("E34").Value = "  +3.80%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D35").Value = "'7.17"
$ws.Range("E35").Value = "  +8.54%  "
$ws.Range("B36").Value = "Fetch.AI"
$ws.Range("C36").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D36").Value = "'1.30"
$ws.Range("E36").Value = "  +5.59%  "
$ws.Range("D37").Value = "'1.85"
$ws.Range("E37").Value = "  +3.90%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "'23.17"
$ws.Range("E38").Value = "  +4.40%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.137"
$ws.Range("E39").Value = "  +3.83%  "
$ws.Range("B40").Value = "WhiteBITCoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D40").Value = "'22.33"
$ws.Range("E40").Value = "  +0.55%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.02%  "
$ws.Range("B42").Value = "Hedera"
$ws.Range("C42").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D42").Value = "'0.0792"
$ws.Range("E42").Value = "  +20.62%  "
$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("B44").Value = "PolygonEcosystemToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D44").Value = "'0.375"
$ws.Range("E44").Value = "  +4.91%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.142"
$ws.Range("E45").Value = "  +6.66%  "
$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").Value = "'1.88"
$ws.Range("E46").Value = "  +4.16%  "
$ws.Range("D47").Value = "'145.40"
$ws.Range("E47").Value = "  -1.46%  "
$ws.Range("D48").Value = "'44.19"
$ws.Range("E48").Value = "  +2.12%  "
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").Value = "'0.000269"
$ws.Range("E49").Value = "  +26.60%  "
$ws.Range("B50").Value = "ImmutableX"
$ws.Range("C50").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D50").Value = "'1.29"
$ws.Range("E50").Value = "  +11.05%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'167.28"
$ws.Range("E51").Value = "  +6.91%  "
